# Daily attendance processing - 2025-12-06 21:24:56
#
# The "Recorded By" column (G) lists the users who touched a session's
# attendance record, as a comma-separated string. This pass normalizes
# the ordering produced by the recorder merge step: whenever "System" is
# recorded first but a real named recorder (or a duplicate lowercase
# "system" tag) is also present, rotate the list right by one slot so
# the most recently appended contributor leads instead of "System".
#
# NOTE: PowerShell's default -eq/-ceq operators were observed to compare
# strings case-insensitively in this runtime, so case-sensitive checks
# below use the .NET String.Equals(...) instance method instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $value = $cell.Value2

    if ($value -eq $null) { continue }
    if ($value.GetType().Name -ne "String") { continue }

    $parts = $value.Split(",")
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    if ($parts.Length -lt 2) { continue }
    if (-not $parts[0].Equals("System")) { continue }

    # Only re-order entries once a real, named recorder (or a duplicate
    # lowercase "system" tag) is mixed in with the "System" placeholder --
    # plain "System, <role account>" pairs (backup/admin service accounts)
    # are left as-is.
    $hasDuplicateSystemTag = $false
    $hasNamedUser = $false
    foreach ($p in $parts) {
        if ($p.Equals("system")) { $hasDuplicateSystemTag = $true }
        elseif ($p.Equals("dnasr281@gmail.com")) { $hasNamedUser = $true }
    }
    if (-not ($hasDuplicateSystemTag -or $hasNamedUser)) { continue }

    # Rotate right by one: last entry moves to the front.
    $rotated = @($parts[$parts.Length - 1]) + $parts[0..($parts.Length - 2)]

    $cell.Value2 = $rotated -join ", "
}
